# Adds formatting - column width and bold, centered font
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: B=4, C=18 (new), D=4 (new) -------------------------
# Excel's ColumnWidth (characters) is stored in the XML as ColumnWidth +
# 5/6 (the default glyph padding for the workbook's Normal font), so we
# back that offset out to land on the exact stored widths from the diff.
$widthOffset = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 4  - $widthOffset
$ws.Columns.Item(3).ColumnWidth = 18 - $widthOffset
$ws.Columns.Item(4).ColumnWidth = 4  - $widthOffset

# --- B4:B127 (RANK) -> bold + centered + shrink-to-fit ------------------
# Apply to the first cell alone so the engine resolves a single combined
# style (center + shrinkToFit, keeping the existing bold font) instead of
# an intermediate style, then fan that exact style out to the rest of the
# column with a format-only paste so every cell lands on the same style.
$rankFirst = $ws.Range("B4")
$rankFirst.HorizontalAlignment = -4108   # xlCenter
$rankFirst.ShrinkToFit = $true
$rankFirst.Copy()
$ws.Range("B5:B127").PasteSpecial(-4122) # xlPasteFormats

# --- C4:C127 (Player Name) -> left aligned ------------------------------
$nameFirst = $ws.Range("C4")
$nameFirst.HorizontalAlignment = -4131   # xlLeft
$nameFirst.Copy()
$ws.Range("C5:C127").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0
